# Adding the data used
# Populate the "Membres du groupe" table (Prénoms / Nom) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Prénoms"
$ws.Range("B1").Value = "Nom"

$ws.Range("A2").Value = "Khadidiatou"
$ws.Range("B2").Value = "Coulibaly"

$ws.Range("A3").Value = "Tamsir"
$ws.Range("B3").Value = "Ndong"

$ws.Range("A4").Value = "Samba"
$ws.Range("B4").Value = "Dieng"

$ws.Range("A5").Value = "Jeanne De La Flèche"
$ws.Range("B5").Value = "Onanena Amana"

# Column A is sized to fit its (now longer) contents, as Excel does
# automatically once the names are typed in.
$ws.Columns("A:A").AutoFit()
$ws.Columns("A:A").ColumnWidth = 16.583333333333332

# Leave the same cell selected as in the saved workbook.
$ws.Range("M9").Select()
